$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet - Latest HO Xliff Generate Date column (G)
$overview.Range("G2").Value = "2016-09-07 08:23:34"
$overview.Range("G3").Value = "2016-09-07 08:23:34"

# zh-cn sheet - Status (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H2").Value = "2016-09-07 08:23:29"
$zhcn.Range("H3").Value = "2016-09-07 08:23:29"
$zhcn.Range("K2").Value = "2016-09-07 08:23:49"
$zhcn.Range("K3").Value = "2016-09-07 08:23:49"

# de-de sheet - Status (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$dede.Range("E2").Value = "mt"
$dede.Range("E3").Value = "mt"
$dede.Range("H2").Value = "2016-09-07 08:23:34"
$dede.Range("H3").Value = "2016-09-07 08:23:34"
$dede.Range("K2").Value = "2016-09-07 08:24:00"
$dede.Range("K3").Value = "2016-09-07 08:24:00"
